$wb = $excel.ActiveWorkbook

# --- Status text update: "Ready for handoff" -> "In Translation" -----------
# "Overview" sheet keeps one status column per locale (zh-cn -> E, de-de -> F)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# Per-locale detail sheets each keep their own "Status" column (column C)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Column width update for the (now narrower) status columns -------------
# The shorter "In Translation" label lets the status columns shrink from
# ~17.22 chars down to ~13.41 chars (this is what a real "AutoFit" /ColumnWidth
# resize produced upstream). The host only persists ColumnWidth at whole-pixel
# granularity (stored width = round(charWidth*6 + 5) / 6), so feed it the
# character width whose pixel bucket lands nearest the 13.4101845877511
# target (the 13.3333.. bucket, centred at charWidth = 12.5) rather than the
# raw target value itself.
$newStatusWidth = 12.5

$wsOverview.Columns("E:F").ColumnWidth = $newStatusWidth
$wsZhCn.Columns("C:C").ColumnWidth = $newStatusWidth
$wsDeDe.Columns("C:C").ColumnWidth = $newStatusWidth
